# Add a new "localdb" macro command category to the "#system" reference
# sheet (used by the data-validation dropdowns on the "macros" sheet), and
# register it in the named ranges so it is usable / autocompleted.
#
# The "#system" sheet keeps one lookup-list per macro category in its own
# column (e.g. column N held "macro" with header + 3 items). We need to
# make room for the brand-new "localdb" category at column N, which pushes
# every following category one column to the right (N->O, O->P, ... AC->AD).
# Column A ("target") is the master, alphabetically sorted list of every
# category name; "localdb" needs to be spliced into it between "json" and
# "macro" (row 14), pushing everything below down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- number of populated rows (including the header in row 1) for each of
# the existing category columns N..AC, keyed by 1-based column index -------
$rowCounts = @{
    14 = 4;    # N  macro
    15 = 2;    # O  mail
    16 = 16;   # P  number
    17 = 16;   # Q  pdf
    18 = 7;    # R  rdbms
    19 = 10;   # S  redis
    20 = 2;    # T  sms
    21 = 5;    # U  sound
    22 = 9;    # V  ssh
    23 = 4;    # W  step
    24 = 127;  # X  web
    25 = 8;    # Y  webalert
    26 = 8;    # Z  webcookie
    27 = 17;   # AA ws
    28 = 8;    # AB ws.async
    29 = 21;   # AC xml
}

# --- 1) shift columns N..AC one column to the right (AC->AD first, so we
# never clobber a column before it has been read) --------------------------
for ($col = 29; $col -ge 14; $col--) {
    $n = $rowCounts[$col]
    for ($r = 1; $r -le $n; $r++) {
        $v = $ws.Cells.Item($r, $col).Value()
        $ws.Cells.Item($r, $col + 1).Value = $v
    }
}

# --- 2) populate the now-vacant column N with the new "localdb" category --
$ws.Cells.Item(1, 14).Value = "localdb"
$ws.Cells.Item(2, 14).Value = "cloneTable(var,source,target)"
$ws.Cells.Item(3, 14).Value = "dropTables(var,tables)"
$ws.Cells.Item(4, 14).Value = "exportCSV(sql,output)"
$ws.Cells.Item(5, 14).Value = "importRecords(var,sourceDb,sql,table)"
$ws.Cells.Item(6, 14).Value = "purge(var)"
$ws.Cells.Item(7, 14).Value = "runSQLs(var,sqls)"

# --- 3) splice "localdb" into the master "target" list in column A, right
# before "macro" (row 14), pushing rows 14..29 down to 15..30 --------------
for ($r = 29; $r -ge 14; $r--) {
    $v = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 1).Value = $v
}
$ws.Cells.Item(14, 1).Value = "localdb"

# --- 4) fix up the workbook-level named ranges: every name that pointed at
# one of the shifted columns needs its reference moved one column right,
# "target" needs its row extent bumped by one, and a brand new "localdb"
# name needs to be created -------------------------------------------------
$wb.Names.Item("macro").RefersTo      = '=''#system''!$O$2:$O$4'
$wb.Names.Item("mail").RefersTo       = '=''#system''!$P$2:$P$2'
$wb.Names.Item("number").RefersTo     = '=''#system''!$Q$2:$Q$16'
$wb.Names.Item("pdf").RefersTo        = '=''#system''!$R$2:$R$16'
$wb.Names.Item("rdbms").RefersTo      = '=''#system''!$S$2:$S$7'
$wb.Names.Item("redis").RefersTo      = '=''#system''!$T$2:$T$10'
$wb.Names.Item("sms").RefersTo        = '=''#system''!$U$2:$U$2'
$wb.Names.Item("sound").RefersTo      = '=''#system''!$V$2:$V$5'
$wb.Names.Item("ssh").RefersTo        = '=''#system''!$W$2:$W$9'
$wb.Names.Item("step").RefersTo       = '=''#system''!$X$2:$X$4'
$wb.Names.Item("target").RefersTo     = '=''#system''!$A$2:$A$30'
$wb.Names.Item("web").RefersTo        = '=''#system''!$Y$2:$Y$127'
$wb.Names.Item("webalert").RefersTo   = '=''#system''!$Z$2:$Z$8'
$wb.Names.Item("webcookie").RefersTo  = '=''#system''!$AA$2:$AA$8'
$wb.Names.Item("ws").RefersTo         = '=''#system''!$AB$2:$AB$17'
$wb.Names.Item("ws.async").RefersTo   = '=''#system''!$AC$2:$AC$8'
$wb.Names.Item("xml").RefersTo        = '=''#system''!$AD$2:$AD$21'

$wb.Names.Add("localdb", '=''#system''!$N$2:$N$7')
